$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unneeded blank placeholder rows (previously rows 50-57, 8 rows
# of unused entries), shifting the totals row (previously row 58) up to row 50.
$ws.Rows("50:57").Delete()

# Row 49 keeps its blank placeholder cells, but it's no longer a real entry row:
# clear out the Rate and the computed Amount formula.
$ws.Range("D49:E49").ClearContents()

# Apply the currency (Amount) number format used elsewhere on the sheet to the
# remaining Amount formula cells E44:E48.
$ws.Range("E44:E48").NumberFormat = "_-""£""* #,##0.00_-;\-""£""* #,##0.00_-;_-""£""* ""-""??_-;_-@_-"

# Update the active selection to reflect the new layout.
$null = $ws.Range("F50").Select()
